$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This change cyclically re-shuffles the species-record data (columns
# A,B,D,E,F,G,H,J,K,L,M,N,Q,R,AF) across rows 53-62, while leaving the
# location/administrative columns (C,I,P,S,T,U,V,W,Y,Z,AA,AB,AD,AE,AG,
# AT,AW,AX,AY) untouched in place. Below we write the full target state
# for every touched row, computed from the unified diff.

# Row 53 (was old row 54's species data)
$ws.Range("A53").Value = 112230603
$ws.Range("B53").Value = 78647
$ws.Range("D53").Value = "LC"
$ws.Range("E53").Value = 6456
$ws.Range("F53").Value = "Skinnlav"
$ws.Range("G53").Value = "Leptogium saturninum"
$ws.Range("H53").Value = "(Dicks.) Nyl."
$ws.Range("Q53").Value = 572018
$ws.Range("R53").Value = 6697738

# Row 54 (was old row 59's species data)
$ws.Range("A54").Value = 112230610
$ws.Range("B54").Value = 90480
$ws.Range("D54").Value = "LC"
$ws.Range("E54").Value = 4769
$ws.Range("F54").Value = "Svavelriska"
$ws.Range("G54").Value = "Lactarius scrobiculatus"
$ws.Range("H54").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q54").Value = 571853
$ws.Range("R54").Value = 6697760

# Row 55 (was old row 53's species data)
$ws.Range("A55").Value = 112230614
$ws.Range("B55").Value = 78647
$ws.Range("D55").Value = "LC"
$ws.Range("E55").Value = 6456
$ws.Range("F55").Value = "Skinnlav"
$ws.Range("G55").Value = "Leptogium saturninum"
$ws.Range("H55").Value = "(Dicks.) Nyl."
$ws.Range("Q55").Value = 571792
$ws.Range("R55").Value = 6697651

# Row 56 (was old row 58's species data); also loses J/K/L/M/N/AF cells
$ws.Range("A56").Value = 112230608
$ws.Range("B56").Value = 99874
$ws.Range("D56").Value = "LC"
$ws.Range("E56").Value = 221235
$ws.Range("F56").Value = "Vårärt"
$ws.Range("G56").Value = "Lathyrus vernus"
$ws.Range("H56").Value = "(L.) Bernh."
$ws.Range("J56").Value = ""
$ws.Range("K56").Value = ""
$ws.Range("L56").Value = ""
$ws.Range("M56").Value = ""
$ws.Range("N56").Value = ""
$ws.Range("AF56").Value = ""
$ws.Range("Q56").Value = 571931
$ws.Range("R56").Value = 6697694

# Row 57 (was old row 62's species data)
$ws.Range("A57").Value = 112230613
$ws.Range("B57").Value = 89553
$ws.Range("D57").Value = "NT"
$ws.Range("E57").Value = 1202
$ws.Range("F57").Value = "Ullticka"
$ws.Range("G57").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H57").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q57").Value = 571799
$ws.Range("R57").Value = 6697620

# Row 58 (was old row 60's species data)
$ws.Range("A58").Value = 112230606
$ws.Range("B58").Value = 56575
$ws.Range("D58").Value = "NT"
$ws.Range("E58").Value = 103021
$ws.Range("F58").Value = "Talltita"
$ws.Range("G58").Value = "Poecile montanus"
$ws.Range("H58").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("Q58").Value = 571961
$ws.Range("R58").Value = 6697705

# Row 59 (was old row 55's species data)
$ws.Range("A59").Value = 112230604
$ws.Range("B59").Value = 102192
$ws.Range("D59").Value = "LC"
$ws.Range("E59").Value = 222412
$ws.Range("F59").Value = "Tibast"
$ws.Range("G59").Value = "Daphne mezereum"
$ws.Range("H59").Value = "L."
$ws.Range("Q59").Value = 571996
$ws.Range("R59").Value = 6697876

# Row 60 (was old row 61's species data)
$ws.Range("A60").Value = 112230611
$ws.Range("B60").Value = 4711
$ws.Range("D60").Value = "LC"
$ws.Range("E60").Value = 100299
$ws.Range("F60").Value = "Thomsons trägnagare"
$ws.Range("G60").Value = "Cacotemnus thomsoni"
$ws.Range("H60").Value = "(Kraatz, 1881)"
$ws.Range("Q60").Value = 571834
$ws.Range("R60").Value = 6697641

# Row 61 (was old row 57's species data)
$ws.Range("A61").Value = 112230605
$ws.Range("B61").Value = 99874
$ws.Range("D61").Value = "LC"
$ws.Range("E61").Value = 221235
$ws.Range("F61").Value = "Vårärt"
$ws.Range("G61").Value = "Lathyrus vernus"
$ws.Range("H61").Value = "(L.) Bernh."
$ws.Range("Q61").Value = 571995
$ws.Range("R61").Value = 6697876

# Row 62 (was old row 56's species data); gains J/K/L/M/N/AF cells
$ws.Range("A62").Value = 112230612
$ws.Range("B62").Value = 12274
$ws.Range("D62").Value = "NT"
$ws.Range("E62").Value = 102016
$ws.Range("F62").Value = "Gropig brunbagge"
$ws.Range("G62").Value = "Zilora ferruginea"
$ws.Range("H62").Value = "(Paykull, 1798)"
$ws.Range("K62").Value = "larv/nymf"
$ws.Range("Q62").Value = 571800
$ws.Range("R62").Value = 6697623
